$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# About sheet: add "HK Notes" section
# ----------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A10").Value = "HK Notes"
$wsAbout.Range("A10").Font.Bold = $true
$wsAbout.Range("A11").Value = "this needs some work."
$wsAbout.Range("A12").Value = "Due to lack of data sources we simply scale from US values by using the ratio of HK's total capacity to US total capacity."

# ----------------------------------------------------------------------
# Data sheet: add the HK scale-factor block (US/HK capacity ratio), and
# highlight the "New transmission quantity in 2050" result
# ----------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# New block: US total capacity 2017
$wsData.Range("A22").Value = "US total capacity 2017:"
$wsData.Range("B22").Value = 1186943.8999999999
$wsData.Range("B22").NumberFormat = "#,##0.00"
$wsData.Range("B22").WrapText = $true
$wsData.Range("B22").VerticalAlignment = -4108
$wsData.Range("C22").Value = "MW"
$wsData.Range("C22").WrapText = $true
$wsData.Range("C22").VerticalAlignment = -4108

$wsData.Range("A23").Value = "https://www.eia.gov/electricity/annual/html/epa_04_03.html"

# New block: HK total capacity
$wsData.Range("A25").Value = "HK total capacity:"
$wsData.Range("B25").Value = 11780
$wsData.Range("C25").Value = "MW"

$wsData.Range("A26").Value = "(from start year capacity sheet)"

# New block: scale factor = HK total capacity / US total capacity
$wsData.Range("A28").Value = "scale factor:"
$wsData.Range("B28").Formula = '=B25/B22'

# ----------------------------------------------------------------------
# BTC sheet: scale the 2010 transmission capacity by the new HK scale
# factor, and mark the seed cell with the accent fill
# ----------------------------------------------------------------------
$wsBTC = $wb.Worksheets.Item("BTC")
$wsBTC.Range("B2").Formula = '=Data!B12*Data!$B$28'
$wsBTC.Range("B2").Interior.ThemeColor = 7

# Existing "New transmission quantity in 2050" cell gets the same accent fill
$wsData.Range("B17").Interior.ThemeColor = 7

# ----------------------------------------------------------------------
# Restore the cursor/selection on each sheet the way the authors left it,
# finishing back on "About" so that sheet stays the active tab
# ----------------------------------------------------------------------
$wsBTC.Range("C2").Select() | Out-Null
$wsData.Range("B17").Select() | Out-Null
$wsAbout.Range("F18").Select() | Out-Null
